$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '66.941.91'
$ws.Cells.Item(2, 5).Value = '  -8.02%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.638.32'
$ws.Cells.Item(3, 5).Value = '  -8.00%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '573.46'
$ws.Cells.Item(5, 5).Value = '  -5.72%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '168.15'
$ws.Cells.Item(6, 5).Value = '  -1.36%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '3.623.52'
$ws.Cells.Item(7, 5).Value = '  -8.24%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.614'
$ws.Cells.Item(8, 5).Value = '  -10.16%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.52%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.691'
$ws.Cells.Item(10, 5).Value = '  -12.05%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.158'
$ws.Cells.Item(11, 5).Value = '  -13.10%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '50.23'
$ws.Cells.Item(12, 5).Value = '  -11.11%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000282'
$ws.Cells.Item(13, 5).Value = '  -14.00%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '10.27'
$ws.Cells.Item(14, 5).Value = '  -11.30%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '4.238.33'
$ws.Cells.Item(15, 5).Value = '  -7.54%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.676.81'
$ws.Cells.Item(16, 5).Value = '  -7.25%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'TRON'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.126'
$ws.Cells.Item(17, 5).Value = '  -3.47%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'Chainlink'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '19.05'
$ws.Cells.Item(18, 5).Value = '  -9.97%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.65'
$ws.Cells.Item(19, 5).Value = '  -10.01%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '1.10'
$ws.Cells.Item(20, 5).Value = '  -11.12%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '66.894.20'
$ws.Cells.Item(21, 5).Value = '  -8.05%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '399.49'
$ws.Cells.Item(22, 5).Value = '  -10.42%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '4.43'
$ws.Cells.Item(23, 5).Value = '  -8.84%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '86.40'
$ws.Cells.Item(24, 5).Value = '  -10.15%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.98'
$ws.Cells.Item(25, 5).Value = '  -11.04%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '12.52'
$ws.Cells.Item(26, 5).Value = '  -11.59%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.50'
$ws.Cells.Item(27, 5).Value = '  -7.18%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '5.97'
$ws.Cells.Item(28, 5).Value = '  +1.35%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '3.67'
$ws.Cells.Item(29, 5).Value = '  -14.10%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '9.27'
$ws.Cells.Item(30, 5).Value = '  -10.80%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '32.04'
$ws.Cells.Item(31, 5).Value = '  -10.44%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.36'
$ws.Cells.Item(32, 5).Value = '  -8.00%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '12.21'
$ws.Cells.Item(33, 5).Value = '  -11.97%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '64.26'
$ws.Cells.Item(34, 5).Value = '  -7.12%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.114'
$ws.Cells.Item(35, 5).Value = '  -11.09%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '42.33'
$ws.Cells.Item(36, 5).Value = '  -16.33%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '581.54'
$ws.Cells.Item(37, 5).Value = '  -8.68%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0₃0871'
$ws.Cells.Item(38, 5).Value = '  -12.84%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.997'
$ws.Cells.Item(39, 5).Value = '  -0.35%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 5).Value = '  +0.07%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'TheGraph'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.388'
$ws.Cells.Item(41, 5).Value = '  -9.34%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.131'
$ws.Cells.Item(42, 5).Value = '  -10.19%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.93'
$ws.Cells.Item(43, 5).Value = '  -14.97%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.59'
$ws.Cells.Item(44, 5).Value = '  -1.05%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0427'
$ws.Cells.Item(45, 5).Value = '  -10.93%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.81'
$ws.Cells.Item(46, 5).Value = '  -11.60%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '9.00'
$ws.Cells.Item(47, 5).Value = '  -14.72%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.762.24'
$ws.Cells.Item(48, 5).Value = '  -2.45%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.131'
$ws.Cells.Item(49, 5).Value = '  -11.37%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '3.11'
$ws.Cells.Item(50, 5).Value = '  -8.53%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.61'
$ws.Cells.Item(51, 5).Value = '  -12.47%  '
